# "added more test case for excel read"
#
# Adds a batch of new test columns/values to the "test-3" sheet, making it
# the active sheet/tab (it previously was "test-1").

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # "test-3"

# Insert a new (blank) column before the existing column L. This shifts the
# old L/M columns (and their custom widths) one slot to the right, becoming
# M/N, exactly mirroring the target column layout.
$ws3.Columns.Item(12).Insert()

# The freshly inserted column (L) needs its own custom width, matching the
# new column O that is appended further along with fresh test data.
$ws3.Columns.Item(15).ColumnWidth = 11.5

# --- Row 1 (header row) ------------------------------------------------
$ws3.Range("L1").Value = "__EMPTY"
$ws3.Range("O1").Value = "__EMPTY_2"
$ws3.Range("P1").Value = " "
$ws3.Range("Q1").Value = "  "
$ws3.Range("R1").Value = " "

# --- Row 2 (data row) ---------------------------------------------------
$ws3.Range("O2").Value = "E2"
$ws3.Range("P2").Value = "S_1"
$ws3.Range("Q2").Value = "S_2"

$ws3.Range("F2").Value = "e2"
$ws3.Range("G2").Value = "e3"
$ws3.Range("H2").Value = "e4"
$ws3.Range("I2").Value = "e5"
$ws3.Range("J2").Value = "e6"
$ws3.Range("L2").Value = "e8"
$ws3.Range("R2").Value = "S_1_1"

# Make "test-3" the active sheet/tab with E6 selected, matching the new
# workbook-level active-tab and sheet-level selection state.
$ws3.Select()
$ws3.Range("E6").Select()
